# Add the "sein" (to be) verb-conjugation row to the Verbs table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Verbs")

# Write the new row's values in the same left-to-right-ish order the
# original author typed them in (verb, then the conjugation columns,
# then meaning/IPA last) so new shared-string entries land in the same
# order as the authored workbook.
$ws.Range("A4").Value = "sein"
$ws.Range("D4").Value = "bin"
$ws.Range("E4").Value = "bist"
$ws.Range("F4").Value = "ist"
$ws.Range("G4").Value = "sind"
$ws.Range("H4").Value = "seid"
$ws.Range("I4").Value = "sind"
$ws.Range("J4").Value = "sind"
$ws.Range("B4").Value = "to be"
$ws.Range("C4").Value = "zain"

# Match the formatting used by the existing data rows (style applied to
# the conjugation columns D:J).
$ws.Range("D3:J3").Copy()
$ws.Range("D4:J4").PasteSpecial(-4122)

# Grow the Excel Table / AutoFilter range to include the new row.
$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:J4"))

# Keep the workbook-level _FilterDatabase defined name in sync with the
# table's new extent.
foreach ($n in $wb.Names) {
    $n.RefersTo = "=Verbs!`$A`$1:`$J`$4"
}

# Match the author's final selection in the sheet.
$ws.Range("D5").Select()
